# Doing Updates for Financials
# Insert a new column D (new fiscal year 2018-12-31 data) ahead of the
# existing Period Ending / financial data columns, shifting everything
# that was in D:K one column to the right (E:L), then populate the new
# column D with the new year's figures for each of the three statements
# (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before D; existing D:K data (and formatting) shifts to E:L.
$ws.Columns("D").Insert()

# The freshly inserted column D has no explicit formatting yet (default style).
# Copy the number formats/styles from column E (which now holds what used to be
# column D) into the new column D so the new cells match the surrounding table
# (date format for header rows, number format for data rows).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# ---- Income Statement (new FY2018-12-31 column) ----
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 34000
$ws.Range("D9").Value = 600
$ws.Range("D10").Value = 33400
$ws.Range("D12").Value = 39300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 149100
$ws.Range("D18").Value = -115000
$ws.Range("D20").Value = 3100
$ws.Range("D21").Value = -110400
$ws.Range("D22").Value = 19100
$ws.Range("D23").Value = -131000
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -131000
$ws.Range("D27").Value = -131000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3100
$ws.Range("D33").Value = -131000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -131000

# ---- Balance Sheet (new FY2018-12-31 column) ----
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 56600
$ws.Range("D42").Value = 154300
$ws.Range("D43").Value = 6400
$ws.Range("D44").Value = 5100
$ws.Range("D45").Value = 6000
$ws.Range("D46").Value = 228400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 3700
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 234800
$ws.Range("D57").Value = 6600
$ws.Range("D58").Value = 1700
$ws.Range("D59").Value = 16000
$ws.Range("D60").Value = 24300
$ws.Range("D61").Value = 117500
$ws.Range("D62").Value = 3200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 144900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -342700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 89900
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (new FY2018-12-31 column) ----
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -131000
$ws.Range("D83").Value = 1500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -104200
$ws.Range("D91").Value = -1100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -69300
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 138900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -34700
